$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D and E cells we touch remain plain text (avoid Excel numeric auto-conversion)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.327.39'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.099.08'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '388.82'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.82'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.91'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.586.47'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.53'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.104.61'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.431.18'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.43'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.87'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.38'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.24'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.165'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.20%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.43'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '36.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +7.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0477'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.06'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.38'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.290'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '130.93'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.87'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.68'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.49'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.09'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.85%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.073.94'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.940'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +19.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0328'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.75%  '
